# "updated scope value plot"
# roi.csv.xlsx holds acquisition ROI definitions (Name, Y1, Y2, X1, X2,
# ImageSizeY, ImageSizeX, Angle, SubRoiCenterSize, SubRoiNRowColumn,
# SubRoiSeparation) keyed one row per scope/config name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 25 ("BMPDloopTof5000"): the ROI box moved/resized - refresh the
# Y/X bounds and the SubRoiCenterSize string, name/grid/separation unchanged.
$ws.Range("B25").Value = 487
$ws.Range("C25").Value = 1545
$ws.Range("D25").Value = 1299
$ws.Range("E25").Value = 1473
$ws.Range("I25").Value = "[939 1345 180 100]"

# Row 30 (new): clone of the refreshed "BMPDloopTof5000" ROI, saved under a
# new name "BMPDloopTof4000" with a smaller SubRoiSeparation.
$ws.Range("A30").Value = "BMPDloopTof4000"
$ws.Range("B30").Value = 487
$ws.Range("C30").Value = 1545
$ws.Range("D30").Value = 1299
$ws.Range("E30").Value = 1473
$ws.Range("F30").Value = 2160
$ws.Range("G30").Value = 2560
$ws.Range("H30").Value = 2.2999999999999998
$ws.Range("I30").Value = "[939 1345 180 100]"
$ws.Range("J30").Value = "[2 1]"
$ws.Range("K30").Value = "[600 100]"
